$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.340.29'
$ws.Range('E2').Value = '  +1.22%  '
$ws.Range('D3').Value = '1.683.98'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').Value = '218.37'
$ws.Range('E5').Value = '  +0.71%  '
$ws.Range('D6').Value = '0.5528'
$ws.Range('E6').Value = '  +8.35%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('D8').Value = '0.2700'
$ws.Range('E8').Value = '  +1.63%  '
$ws.Range('D9').Value = '0.06498'
$ws.Range('E9').Value = '  +1.50%  '
$ws.Range('D10').Value = '22.12'
$ws.Range('E10').Value = '  +1.48%  '
$ws.Range('D11').Value = '0.07563'
$ws.Range('E11').Value = '  +1.68%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.700.99'
$ws.Range('E12').Value = '  +1.88%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.546'
$ws.Range('E13').Value = '  +0.75%  '
$ws.Range('D14').Value = '0.5810'
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').Value = '0.000008455'
$ws.Range('E15').Value = '  -1.31%  '
$ws.Range('D16').Value = '65.20'
$ws.Range('E16').Value = '  +1.44%  '
$ws.Range('D17').Value = '26.369.07'
$ws.Range('E17').Value = '  +1.09%  '
$ws.Range('D18').Value = '4.938'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').Value = '10.92'
$ws.Range('E20').Value = '  +1.36%  '
$ws.Range('D21').Value = '191.41'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = '6.238'
$ws.Range('E22').Value = '  +0.54%  '
$ws.Range('D23').Value = '1.009'
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('D24').Value = '147.77'
$ws.Range('E24').Value = '  +2.12%  '
$ws.Range('D25').Value = '0.1326'
$ws.Range('E25').Value = '  +10.55%  '
$ws.Range('D26').Value = '7.897'
$ws.Range('E26').Value = '  +3.74%  '
$ws.Range('D27').Value = '15.81'
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('D28').Value = '0.06348'
$ws.Range('E28').Value = '  -2.42%  '
$ws.Range('D29').Value = '1.395'
$ws.Range('E29').Value = '  +4.02%  '
$ws.Range('E30').Value = '  +0.42%  '
$ws.Range('D31').Value = '3.593'
$ws.Range('E31').Value = '  +1.54%  '
$ws.Range('D32').Value = '3.582'
$ws.Range('E33').Value = '  +1.24%  '
$ws.Range('D34').Value = '1.042'
$ws.Range('E34').Value = '  +2.21%  '
$ws.Range('D35').Value = '0.6224'
$ws.Range('E35').Value = '  +1.78%  '
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('D37').Value = '2.720'
$ws.Range('E37').Value = '  +1.33%  '
$ws.Range('D38').Value = '6.234'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').Value = '1.113.55'
$ws.Range('E39').Value = '  +2.15%  '
$ws.Range('D40').Value = '0.01629'
$ws.Range('E40').Value = '  +1.72%  '
$ws.Range('D41').Value = '0.8728'
$ws.Range('E41').Value = '  +1.07%  '
$ws.Range('E42').Value = '  +0.69%  '
$ws.Range('D43').Value = '100.72'
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('D44').Value = '1.833.46'
$ws.Range('E44').Value = '  +0.99%  '
$ws.Range('D45').Value = '0.00000000109'
$ws.Range('E45').Value = '  -4.23%  '
$ws.Range('D46').Value = '57.38'
$ws.Range('E46').Value = '  +1.75%  '
$ws.Range('E47').Value = '  +1.96%  '
$ws.Range('E48').Value = '  -0.21%  '
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('D50').Value = '0.4295'
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').Value = '6.079'
$ws.Range('E51').Value = '  -0.15%  '
